$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos ("minuts" -> "minutes") and update the value for row 6 (B6)
# which changes from "8 hours 30 minuts" to "14 hours 30 minutes".
$ws.Range("B3").Value = "4 hours 33 minutes"
$ws.Range("B4").Value = "2 hours 30 minutes"
$ws.Range("B5").Value = "7 hours 30 minutes"
$ws.Range("B6").Value = "14 hours 30 minutes"

# Move the active selection from B7 to F7.
$ws.Range("F7").Select()
